$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.779.52"
$ws.Range("E2").Value = "  +2.28%  "

$ws.Range("D3").Value = "3.812.19"
$ws.Range("E3").Value = "  +0.36%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.28%  "

$ws.Range("D5").Value = "632.61"
$ws.Range("E5").Value = "  +5.24%  "

$ws.Range("D6").Value = "165.20"
$ws.Range("E6").Value = "  -0.06%  "

$ws.Range("D7").Value = "3.809.88"
$ws.Range("E7").Value = "  +0.36%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("E9").Value = "  +0.80%  "

$ws.Range("E10").Value = "  +1.45%  "

$ws.Range("D11").Value = "0.454"
$ws.Range("E11").Value = "  +0.69%  "

$ws.Range("D12").Value = "6.70"
$ws.Range("E12").Value = "  +3.29%  "

$ws.Range("D13").Value = "0.0000250"
$ws.Range("E13").Value = "  +0.13%  "

$ws.Range("D14").Value = "35.94"
$ws.Range("E14").Value = "  +0.37%  "

$ws.Range("D15").Value = "4.464.76"
$ws.Range("E15").Value = "  +0.69%  "

$ws.Range("D16").Value = "3.804.55"
$ws.Range("E16").Value = "  +0.45%  "

$ws.Range("D17").Value = "69.820.65"
$ws.Range("E17").Value = "  +2.41%  "

$ws.Range("D18").Value = "18.08"
$ws.Range("E18").Value = "  -1.96%  "

$ws.Range("D19").Value = "7.15"
$ws.Range("E19").Value = "  +0.94%  "

$ws.Range("E20").Value = "  -0.74%  "

$ws.Range("D21").Value = "469.72"
$ws.Range("E21").Value = "  +1.79%  "

$ws.Range("D22").Value = "9.71"
$ws.Range("E22").Value = "  -0.14%  "

$ws.Range("D23").Value = "0.708"
$ws.Range("E23").Value = "  +1.16%  "

$ws.Range("E24").Value = "  +1.36%  "

$ws.Range("D25").Value = "83.75"
$ws.Range("E25").Value = "  +0.89%  "

$ws.Range("D26").Value = "12.19"
$ws.Range("E26").Value = "  +1.00%  "

$ws.Range("E27").Value = "  +3.17%  "

$ws.Range("D28").Value = "10.09"
$ws.Range("E28").Value = "  +1.02%  "

$ws.Range("E29").Value = "  -0.01%  "

$ws.Range("D30").Value = "3.971.37"
$ws.Range("E30").Value = "  +0.61%  "

$ws.Range("D31").Value = "2.70"
$ws.Range("E31").Value = "  +1.68%  "

$ws.Range("D32").Value = "2.24"
$ws.Range("E32").Value = "  +0.33%  "

$ws.Range("D33").Value = "7.35"
$ws.Range("E33").Value = "  -0.18%  "

$ws.Range("D34").Value = "29.22"
$ws.Range("E34").Value = "  -0.50%  "

$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").Value = "9.07"
$ws.Range("E36").Value = "  +0.64%  "

$ws.Range("B37").Value = "RenzoRestakedETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D37").Value = "3.757.65"
$ws.Range("E37").Value = "  +0.23%  "

$ws.Range("E38").Value = "  +3.48%  "

$ws.Range("E39").Value = "  +9.01%  "

$ws.Range("D40").Value = "3.35"
$ws.Range("E40").Value = "  +1.34%  "

$ws.Range("D41").Value = "5.90"
$ws.Range("E41").Value = "  +1.57%  "

$ws.Range("E42").Value = "  -0.81%  "

$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.12%  "

$ws.Range("D45").Value = "155.38"
$ws.Range("E45").Value = "  +2.38%  "

$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "1.96"
$ws.Range("E46").Value = "  +4.24%  "

$ws.Range("B47").Value = "TheGraph"
$ws.Range("C47").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D47").Value = "0.301"
$ws.Range("E47").Value = "  +0.17%  "

$ws.Range("D48").Value = "47.21"
$ws.Range("E48").Value = "  -0.59%  "

$ws.Range("D49").Value = "43.24"
$ws.Range("E49").Value = "  -0.62%  "

$ws.Range("E50").Value = "  +3.57%  "

$ws.Range("D51").Value = "8.46"
$ws.Range("E51").Value = "  +1.20%  "
